$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.506.89'
$ws.Range("E2").Value = '  +0.57%  '
$ws.Range("D3").Value = '2.692.24'
$ws.Range("E3").Value = '  +1.86%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.54'
$ws.Range("E5").Value = '  +0.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '160.08'
$ws.Range("E6").Value = '  +2.72%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.543'
$ws.Range("E8").Value = '  +0.28%  '
$ws.Range("D9").Value = '2.691.65'
$ws.Range("E9").Value = '  +1.89%  '
$ws.Range("E10").Value = '  -0.72%  '
$ws.Range("E11").Value = '  -0.54%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.29'
$ws.Range("E13").Value = '  +2.77%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.22'
$ws.Range("E14").Value = '  +0.89%  '
$ws.Range("D15").Value = '3.183.42'
$ws.Range("E15").Value = '  +1.89%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000188'
$ws.Range("E16").Value = '  -0.55%  '
$ws.Range("D17").Value = '68.486.40'
$ws.Range("E17").Value = '  +0.67%  '
$ws.Range("D18").Value = '2.699.86'
$ws.Range("E18").Value = '  +1.99%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.86'
$ws.Range("E19").Value = '  +4.52%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '366.38'
$ws.Range("E20").Value = '  +0.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.63'
$ws.Range("E21").Value = '  +3.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.54'
$ws.Range("E22").Value = '  +2.75%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.88'
$ws.Range("E23").Value = '  +2.14%  '
$ws.Range("E24").Value = '  +2.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.51'
$ws.Range("E25").Value = '  -0.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.02'
$ws.Range("E27").Value = '  +3.43%  '
$ws.Range("D28").Value = '2.829.78'
$ws.Range("E28").Value = '  +1.98%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000104'
$ws.Range("E29").Value = '  +0.44%  '
$ws.Range("E30").Value = '  +0.09%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '571.65'
$ws.Range("E31").Value = '  +3.06%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.22'
$ws.Range("E32").Value = '  +2.65%  '
$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.45'
$ws.Range("E33").Value = '  +3.60%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.94'
$ws.Range("E34").Value = '  +5.49%  '
$ws.Range("E35").Value = '  +2.72%  '
$ws.Range("E36").Value = '  +6.75%  '
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '161.78'
$ws.Range("E38").Value = '  +0.38%  '
$ws.Range("B39").Value = 'EthereumClassic'
$ws.Range("C39").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.87'
$ws.Range("E39").Value = '  +2.34%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.378'
$ws.Range("E40").Value = '  +1.78%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.91'
$ws.Range("E41").Value = '  +2.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.40'
$ws.Range("E42").Value = '  +1.68%  '
$ws.Range("B43").Value = 'WhiteBITCoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.88'
$ws.Range("E43").Value = '  +0.48%  '
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.65'
$ws.Range("E44").Value = '  +1.75%  '
$ws.Range("E46").Value = '  -6.15%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '157.39'
$ws.Range("E47").Value = '  -1.13%  '
$ws.Range("E48").Value = '  +7.39%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.77'
$ws.Range("E49").Value = '  +4.81%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.599'
$ws.Range("E50").Value = '  +6.67%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.99'
$ws.Range("E51").Value = '  -0.03%  '
